# Apply the commit: "code refactoring and loan accounting and charges added"
$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)   # ProductLoan_Input
$ws2 = $wb.Worksheets.Item(2)   # ProductLoan_Output

# --- ProductLoan_Input sheet ---

# Rename the product name (shared string referenced by B1) and shortname text value.
$ws1.Cells.Item(1,2).Value = "438-RBI-EI-DB-SAR-REC-NON-RNI-CTPD-DL-MD-TR-1-EarlyRePayment"

# shortname: was text "kar8" -> now numeric 438
$ws1.Cells.Item(3,2).Value = 438

# nominalinterestratedefault: was 12 -> now 1
$ws1.Cells.Item(11,2).Value = 1

# Insert a new row 22 for preclosureinterestcalculationrule, taking on the
# formatting of the row above it (row 21).
$ws1.Rows.Item(22).Insert()
$ws1.Cells.Item(22,1).Value = "preclosureinterestcalculationrule"
$ws1.Cells.Item(22,2).Value = "Calculate till pre closure date"

# Append the loan accounting / charges rows (31-42), copying the style of an
# existing plain row (row 21: col A style 3, col B style 1) onto the new block
# before filling in values.
$ws1.Range("A21:B21").Copy()
$ws1.Range("A31:B42").PasteSpecial(-4122)

# Column B (account values) was populated first...
$ws1.Cells.Item(31,2).Value = "Cash"
$ws1.Cells.Item(32,2).Value = "Loan portfolio "
$ws1.Cells.Item(33,2).Value = "Interest Receivable "
$ws1.Cells.Item(34,2).Value = "Penalties Receivable "
$ws1.Cells.Item(35,2).Value = "Transfer in Suspence "
$ws1.Cells.Item(36,2).Value = "Fees Receivable"
$ws1.Cells.Item(37,2).Value = "Income from interest"
$ws1.Cells.Item(38,2).Value = "Income from penalties"
$ws1.Cells.Item(39,2).Value = "Income from fees"
$ws1.Cells.Item(40,2).Value = "Income from recovery repayments"
$ws1.Cells.Item(41,2).Value = "Losses Writtenoff "
$ws1.Cells.Item(42,2).Value = "Overpayment Liability"

# ...then column A (account labels) was populated.
$ws1.Cells.Item(31,1).Value = "fundsource"
$ws1.Cells.Item(32,1).Value = "loanprotfolio"
$ws1.Cells.Item(33,1).Value = "interestreceivable"
$ws1.Cells.Item(34,1).Value = "penaltiesreceivable"
$ws1.Cells.Item(35,1).Value = "transferinsuspense"
$ws1.Cells.Item(36,1).Value = "feesreceivable"
$ws1.Cells.Item(37,1).Value = "incomefrominterest"
$ws1.Cells.Item(38,1).Value = "incomefrompenalties"
$ws1.Cells.Item(39,1).Value = "incomefromfees"
$ws1.Cells.Item(40,1).Value = "incomefromrecoveryrepayments"
$ws1.Cells.Item(41,1).Value = "loseswrittenoff"
$ws1.Cells.Item(42,1).Value = "overpaymentliability"

# Update sheet view state to match the edited document (scrolled down a bit,
# selection moved to A44 which is just below the new data).
$ws1.Application.ActiveWindow.ScrollRow = 19
$ws1.Range("A44").Select()

# --- ProductLoan_Output sheet ---
# Product name here mirrors the input sheet, so update it the same way.
$ws2.Cells.Item(1,2).Value = "438-RBI-EI-DB-SAR-REC-NON-RNI-CTPD-DL-MD-TR-1-EarlyRePayment"

# Selection moved down to B29 on this sheet.
$ws2.Range("B29").Select()
